$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values: A1 -> "TestingBLZ123", B1 -> "Theend@1"
$ws.Range("A1").Value = "TestingBLZ123"
$ws.Range("B1").Value = "Theend@1"
